$d = $word.ActiveDocument

# --- 1. Title: "Lead Mailer AI Automation" -> "Weather Report AI Automation"
$d.Content.Find.Execute("Lead Mailer ", $true, $false, $false, $false, $false, `
                         $true, 1, $false, "Weather Report ", 2) | Out-Null

# --- 2. Bold the three "How It Works" list items (first occurrences only):
#        "Form Submission", "Switch Node", "Set Nodes for Coordinates"
$targets = @("Form Submission", "Switch Node", "Set Nodes for Coordinates")
foreach ($label in $targets) {
    for ($i = 1; $i -le $d.Paragraphs.Count; $i++) {
        $para = $d.Paragraphs($i)
        if ($para.Range.Text -eq ($label + "`r")) {
            $para.Range.Bold = $true
            $para.Range.Font.BoldBi = $true
            break
        }
    }
}
